# Daily attendance processing - 2025-12-30 10:34:50
#
# Updates the "Session Analysis Results" sheet:
#  - Recorded Sessions / Total Sessions metrics (L6, L8) and their derived
#    percentages (L9, L10)
#  - Per-group summary rows (16,17,18,24,25,26): Recorded/Pending counts
#    (O,Q) and Coverage %/Avg Attendance % (R,S) for group B1-10
#  - Six individual session rows (27, 47, 67, 182, 202, 222) flip from
#    "Pending" (unrecorded) to "Recorded", with attendance now logged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write literal text into a cell without Excel's autoformatting
# re-interpreting strings that look like percentages/dates/numbers as
# actual numeric values. Doing it via a text formula and then collapsing
# the formula down to its cached value with a values-only paste keeps the
# cell's existing style (number format, fill, font) untouched.
function Set-LiteralText {
    param(
        [string]$CellAddr,
        [string]$Text
    )
    $cell = $ws.Range($CellAddr)
    $escaped = $Text.Replace('"', '""')
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
}

# Helper: re-stamp a whole row's A:I formatting (fill/font) from a known
# "Recorded" template row, without touching the cell contents already
# there -- mirrors the green "Recorded" styling used elsewhere in the sheet.
function Set-RecordedRowFormat {
    param([int]$Row)
    $src = $ws.Range("A2:I2")
    $dst = $ws.Range("A" + $Row + ":I" + $Row)
    $src.Copy()
    $dst.PasteSpecial(-4122)  # xlPasteFormats
    $ws.Application.CutCopyMode = $false
}

# --- Class statistics summary (K/L columns) ---------------------------
$ws.Range("L6").Value = 114          # Recorded Sessions
$ws.Range("L8").Value = 114          # (Total Sessions-derived metric)
Set-LiteralText "L9" "48.7%"
Set-LiteralText "L10" "78.0%"

# --- Per-group table rows for B1-10 (rows 16,17,18,24,25,26) ----------
$groupRows = @(
    @{ Row = 16; S = "71.9%" },
    @{ Row = 17; S = "60.0%" },
    @{ Row = 18; S = "82.9%" },
    @{ Row = 24; S = "71.9%" },
    @{ Row = 25; S = "74.8%" },
    @{ Row = 26; S = "71.7%" }
)

foreach ($entry in $groupRows) {
    $r = $entry.Row
    $ws.Range("O$r").Value = 10
    $ws.Range("Q$r").Value = 10
    Set-LiteralText "R$r" "50.0%"
    Set-LiteralText "S$r" $entry.S
}

# --- Individual sessions that moved from Pending -> Recorded ----------
$sessionRows = @(
    @{ Row = 27;  Students = "27/31" },
    @{ Row = 47;  Students = "13/18" },
    @{ Row = 67;  Students = "15/21" },
    @{ Row = 182; Students = "22/27" },
    @{ Row = 202; Students = "21/29" },
    @{ Row = 222; Students = "19/29" }
)

foreach ($entry in $sessionRows) {
    $r = $entry.Row
    Set-RecordedRowFormat $r
    Set-LiteralText "G$r" "dnasr281@gmail.com"
    Set-LiteralText "H$r" $entry.Students
    Set-LiteralText "I$r" "Recorded"
}
